$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (no numeric auto-conversion) for D-column price cells
# whose new values are valid numeric literals but must stay literal text,
# matching the source feeds plain-string formatting. We temporarily mark
# the cell as Text, write the value, then clear the format again so the
# cell keeps its original (default) style.
$textForceRows = @(5, 6, 7, 8, 9, 10, 11, 13, 14, 15, 18, 20, 21, 22, 23, 25, 26, 27, 28, 29, 30, 32, 33, 34, 35, 36, 37, 38, 39, 40, 41, 42, 43, 44, 46, 48, 49, 50, 51)
foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.555.20"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "1.658.07"
$ws.Range("E3").Value = "  +2.99%  "

$ws.Range("E4").Value = "  -0.19%  "

$ws.Range("D5").Value = "0.9998"
$ws.Range("E5").Value = "  -0.24%  "

$ws.Range("D6").Value = "302.56"
$ws.Range("E6").Value = "  +0.22%  "

$ws.Range("D7").Value = "0.3826"
$ws.Range("E7").Value = "  +1.30%  "

$ws.Range("D8").Value = "0.3601"
$ws.Range("E8").Value = "  +2.23%  "

$ws.Range("D9").Value = "51.11"
$ws.Range("E9").Value = "  -1.71%  "

$ws.Range("D10").Value = "0.08206"
$ws.Range("E10").Value = "  +1.70%  "

$ws.Range("D11").Value = "1.239"
$ws.Range("E11").Value = "  +3.46%  "

$ws.Range("E12").Value = "  -0.41%  "

$ws.Range("D13").Value = "22.46"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "6.487"
$ws.Range("E14").Value = "  +2.20%  "

$ws.Range("D15").Value = "7.509"
$ws.Range("E15").Value = "  +3.73%  "

$ws.Range("E16").Value = "  +1.36%  "

$ws.Range("D17").Value = "1.649.08"
$ws.Range("E17").Value = "  +3.70%  "

$ws.Range("D18").Value = "97.62"
$ws.Range("E18").Value = "  +3.90%  "

$ws.Range("E19").Value = "  +1.20%  "

$ws.Range("D20").Value = "6.814"
$ws.Range("E20").Value = "  +5.63%  "

$ws.Range("D21").Value = "17.66"
$ws.Range("E21").Value = "  +3.06%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.29%  "

$ws.Range("D23").Value = "12.70"
$ws.Range("E23").Value = "  +3.64%  "

$ws.Range("D24").Value = "23.574.62"
$ws.Range("E24").Value = "  +1.41%  "

$ws.Range("D25").Value = "2.522"
$ws.Range("E25").Value = "  -0.63%  "

$ws.Range("D26").Value = "3.043"
$ws.Range("E26").Value = "  -0.85%  "

$ws.Range("D27").Value = "21.22"
$ws.Range("E27").Value = "  +2.07%  "

$ws.Range("D28").Value = "152.63"
$ws.Range("E28").Value = "  +1.01%  "

$ws.Range("D29").Value = "5.247"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").Value = "134.11"
$ws.Range("E30").Value = "  +1.82%  "

$ws.Range("D31").Value = "1.837.36"
$ws.Range("E31").Value = "  +3.58%  "

$ws.Range("D32").Value = "7.170"
$ws.Range("E32").Value = "  +11.52%  "

$ws.Range("D33").Value = "2.248"
$ws.Range("E33").Value = "  +7.20%  "

$ws.Range("D34").Value = "12.08"
$ws.Range("E34").Value = "  +6.67%  "

$ws.Range("D35").Value = "1.060"
$ws.Range("E35").Value = "  -0.15%  "

$ws.Range("D36").Value = "0.02805"
$ws.Range("E36").Value = "  +3.93%  "

$ws.Range("D37").Value = "6.137"
$ws.Range("E37").Value = "  +5.59%  "

$ws.Range("D38").Value = "0.2500"
$ws.Range("E38").Value = "  +2.31%  "

$ws.Range("D39").Value = "0.08787"
$ws.Range("E39").Value = "  +1.37%  "

$ws.Range("D40").Value = "0.07013"
$ws.Range("E40").Value = "  +1.67%  "

$ws.Range("D41").Value = "13.22"
$ws.Range("E41").Value = "  +10.77%  "

$ws.Range("D42").Value = "0.7000"
$ws.Range("E42").Value = "  +2.46%  "

$ws.Range("D43").Value = "1.338"
$ws.Range("E43").Value = "  +1.89%  "

$ws.Range("D44").Value = "15.97"
$ws.Range("E44").Value = "  +4.94%  "

$ws.Range("E45").Value = "  +4.16%  "

$ws.Range("D46").Value = "0.9999"
$ws.Range("E46").Value = "  -0.25%  "

$ws.Range("E47").Value = "  +2.81%  "

$ws.Range("D48").Value = "3.957"
$ws.Range("E48").Value = "  +0.63%  "

$ws.Range("D49").Value = "0.07917"
$ws.Range("E49").Value = "  +0.72%  "

$ws.Range("D50").Value = "128.37"
$ws.Range("E50").Value = "  +0.42%  "

$ws.Range("D51").Value = "1.191"
$ws.Range("E51").Value = "  +2.46%  "

foreach ($r in $textForceRows) {
    $ws.Cells.Item($r, 4).ClearFormats()
}